$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IntakeSheet")

# Data-driven rewrite of the variable table (rows 3-41) per the updated docgen variable list.
# Each entry: row, GroupHeader(A), VariableName(B), Type(D), Description(E)
$rowsData = @(
    @{ Row = 3; A = $null; B = "XX25"; D = "string"; E = "caption number, don't include CI/JV etc. just `"25-101`"" }
    @{ Row = 4; A = $null; B = "affirmdefplrl"; D = "string"; E = "Type the word `"defense`" or `"defenses`"" }
    @{ Row = 5; A = $null; B = "jurisdiction"; D = "string"; E = "County/District" }
    @{ Row = 6; A = $null; B = "servmethod"; D = "string"; E = "Things like e-mail, personal service, e-service." }
    @{ Row = 7; A = $null; B = "venue"; D = "string"; E = "the county name, e.g. `"Garfield County`"" }
    @{ Row = 8; A = $null; B = "docket"; D = "string"; E = "Whether CI, CR, JV, PR, etc." }
    @{ Row = 9; A = "Client"; B = $null; D = $null; E = $null }
    @{ Row = 10; A = $null; B = "birthdate"; D = "date"; E = "Client's birthday" }
    @{ Row = 11; A = $null; B = "clientname"; D = "string"; E = "User-defined: clientname" }
    @{ Row = 12; A = $null; B = "clients"; D = "string"; E = "In letters, pleadings, etc. you would say my clients or my client; type s for multiple clients." }
    @{ Row = 13; A = $null; B = "defendantplural"; D = "string"; E = "Pluralizes defendant vs defendants; type s or leave empty." }
    @{ Row = 14; A = $null; B = "firstname"; D = "string"; E = "Client's first name (legal)" }
    @{ Row = 15; A = $null; B = "firstname_spouse"; D = "string"; E = "User-defined: firstname_spouse" }
    @{ Row = 16; A = $null; B = "lastname"; D = "string"; E = "Client's last name (legal)" }
    @{ Row = 17; A = $null; B = "lastname_spouse"; D = "string"; E = "User-defined: lastname_spouse" }
    @{ Row = 18; A = $null; B = "plaintiffpossessive"; D = "string"; E = "Type Plaintiff's or Plaintiffs' for possesive." }
    @{ Row = 19; A = $null; B = "plurals"; D = "string"; E = "Pluralizes words following defendant, e.g. Defendant jumps or Defendants jump; type s or leave blank." }
    @{ Row = 20; A = $null; B = "spousefirstname"; D = "string"; E = "Lauren" }
    @{ Row = 21; A = $null; B = "spouselastname"; D = "string"; E = "d" }
    @{ Row = 22; A = "Derived"; B = $null; D = $null; E = $null }
    @{ Row = 23; A = $null; B = "clientname_spouse"; D = "string"; E = "User-defined: clientname_spouse" }
    @{ Row = 24; A = $null; B = "defendantcaption"; D = "string"; E = "Derived: wholename_client and wholename_clientspouse" }
    @{ Row = 25; A = $null; B = "defendantscaption"; D = "string"; E = "User-defined: defendantscaption" }
    @{ Row = 26; A = $null; B = "wholename_client"; D = "string"; E = "Derived: firstname  lastname" }
    @{ Row = 27; A = $null; B = "wholename_clientspouse"; D = "string"; E = "Derived: spousefirstname  spouselastname" }
    @{ Row = 28; A = "General"; B = $null; D = $null; E = $null }
    @{ Row = 29; A = $null; B = "Jurisdiction"; D = "string"; E = "User-defined: Jurisdiction" }
    @{ Row = 30; A = $null; B = "Jurisdiction_upper"; D = "string"; E = "User-defined: Jurisdiction_upper" }
    @{ Row = 31; A = $null; B = "clientname_combo"; D = "string"; E = "User-defined: clientname_combo" }
    @{ Row = 32; A = $null; B = "jurisdiction_upper"; D = "string"; E = "User-defined: jurisdiction_upper" }
    @{ Row = 33; A = $null; B = "plaintiffname"; D = "string"; E = "User-defined: plaintiffname" }
    @{ Row = 34; A = $null; B = "plaintiffplural"; D = "string"; E = "User-defined: plaintiffplural" }
    @{ Row = 35; A = $null; B = "venue_upper"; D = "string"; E = "User-defined: venue_upper" }
    @{ Row = 36; A = "Internal"; B = $null; D = $null; E = $null }
    @{ Row = 37; A = $null; B = "matterdesc"; D = "string"; E = "Basic description of case/legal service" }
    @{ Row = 38; A = $null; B = "matterid"; D = "string"; E = "Account/Clio No." }
    @{ Row = 39; A = "Opposing Party"; B = $null; D = $null; E = $null }
    @{ Row = 40; A = $null; B = "CCID"; D = "string"; E = "Client's creditor number." }
    @{ Row = 41; A = $null; B = $null; D = $null; E = $null }
)

foreach ($item in $rowsData) {
    $r = $item.Row
    if ($item.A -ne $null) { $ws.Cells.Item($r, 1).Value = $item.A } else { $ws.Cells.Item($r, 1).Value = "" }
    if ($item.B -ne $null) { $ws.Cells.Item($r, 2).Value = $item.B } else { $ws.Cells.Item($r, 2).Value = "" }
    if ($item.D -ne $null) { $ws.Cells.Item($r, 4).Value = $item.D } else { $ws.Cells.Item($r, 4).Value = "" }
    if ($item.E -ne $null) { $ws.Cells.Item($r, 5).Value = $item.E } else { $ws.Cells.Item($r, 5).Value = "" }
}
